$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (existing "Geanta Mickey Mouse" row) - fill in the Image/Link/Price that were blank
$ws.Range("B5").Value = "https://www.desigual.com/dw/image/v2/BCVV_PRD/on/demandware.static/-/Sites-desigual-m-catalog/default/dwfb405b16/images/B2C/23SAXP44_2000_3.jpg?sfrm=jpg&v=webp10&sw=375"
$ws.Range("C5").Value = "https://www.desigual.com/de_CH/23SAXP442000.html"
$ws.Range("D5").Value = "85.95 CHF"

# Row 6 (new row) - "Perie par (culoare COCONUT)"
$ws.Range("C6").Value = "https://wetbrush.com/collections/go-green/products/go-green-treatment-and-shine-brush?variant=44553084469494"
$ws.Range("B6").Value = "https://wetbrush.com/cdn/shop/products/GO_GREEN_TREATMENT-Paddle-ORANGE-Hair_Brush-BIO833COCON-Wet_Brush-Front_2048x2048.jpg?v=1667403413"
$ws.Range("A6").Value = "Perie par (culoare COCONUT)"
$ws.Range("D6").Value = "16.99 USD"

# Update the active selection to match the author's saved cursor position
$ws.Range("C14").Select()
